# Append two new data rows (118 and 119) to the end of the sheet,
# matching the structure of the existing rows (date/volume/high/low/open/close/adj_close/ticker).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 118
$ws.Cells.Item(118, 1).Value = 45454.2916666667
$ws.Cells.Item(118, 2).Value = 157000
$ws.Cells.Item(118, 3).Value = 3.40000009536743
$ws.Cells.Item(118, 4).Value = 3.07999992370605
$ws.Cells.Item(118, 5).Value = 3.16000008583069
$ws.Cells.Item(118, 6).Value = 3.1800000667572
$ws.Cells.Item(118, 7).Value = "3.1800000667572"
$ws.Cells.Item(118, 8).Value = "AGAIN.MI"

# Row 119
$ws.Cells.Item(119, 1).Value = 45455.2916666667
$ws.Cells.Item(119, 2).Value = 0
$ws.Cells.Item(119, 3).Value = 3.1800000667572
$ws.Cells.Item(119, 4).Value = 3.1800000667572
$ws.Cells.Item(119, 5).Value = 3.1800000667572
$ws.Cells.Item(119, 6).Value = 3.1800000667572
$ws.Cells.Item(119, 7).Value = "3.1800000667572"
$ws.Cells.Item(119, 8).Value = "AGAIN.MI"

# Match the date-time style used by the rest of column A (reuse existing style, xlPasteFormats = -4122)
$ws.Cells.Item(117, 1).Copy()
$ws.Cells.Item(118, 1).PasteSpecial(-4122)
$ws.Cells.Item(119, 1).PasteSpecial(-4122)
